# Fruta / hortaliza, semanal
# Inserts two new weekly price records (rows) into the Arándano (blue)
# dataset for "Vega Central Mapocho de Santiago", pushing the existing
# rows down to make room - first insertion before the current row 87,
# second insertion before the current row 94.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert #1: new record dated 2021-09-10 (serial 44449), before row 87 ---
$ws.Rows.Item(87).Insert()

$ws.Cells.Item(87, 1).Value = 9
$ws.Cells.Item(87, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(87, 3).Value = "Metropolitana"
$ws.Cells.Item(87, 4).Value = 44449
$ws.Cells.Item(87, 5).Value = 13
$ws.Cells.Item(87, 6).Value = "Fruta"
$ws.Cells.Item(87, 7).Value = 100101
$ws.Cells.Item(87, 8).Value = "Berries"
$ws.Cells.Item(87, 9).Value = 100101001
$ws.Cells.Item(87, 10).Value = "Arándano (blue)"
$ws.Cells.Item(87, 11).Value = "Sin especificar"
$ws.Cells.Item(87, 12).Value = "Primera"
$ws.Cells.Item(87, 13).Value = 65
$ws.Cells.Item(87, 14).Value = 16000
$ws.Cells.Item(87, 15).Value = 16000
$ws.Cells.Item(87, 16).Value = 16000
$ws.Cells.Item(87, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(87, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(87, 19).Value = 8000
$ws.Cells.Item(87, 20).Value = 2

# --- Insert #2: new record dated 2021-09-09 (serial 44448), before row 94 ---
$ws.Rows.Item(94).Insert()

$ws.Cells.Item(94, 1).Value = 9
$ws.Cells.Item(94, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(94, 3).Value = "Metropolitana"
$ws.Cells.Item(94, 4).Value = 44448
$ws.Cells.Item(94, 5).Value = 13
$ws.Cells.Item(94, 6).Value = "Fruta"
$ws.Cells.Item(94, 7).Value = 100101
$ws.Cells.Item(94, 8).Value = "Berries"
$ws.Cells.Item(94, 9).Value = 100101001
$ws.Cells.Item(94, 10).Value = "Arándano (blue)"
$ws.Cells.Item(94, 11).Value = "Sin especificar"
$ws.Cells.Item(94, 12).Value = "Primera"
$ws.Cells.Item(94, 13).Value = 50
$ws.Cells.Item(94, 14).Value = 16000
$ws.Cells.Item(94, 15).Value = 16000
$ws.Cells.Item(94, 16).Value = 16000
$ws.Cells.Item(94, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(94, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(94, 19).Value = 8000
$ws.Cells.Item(94, 20).Value = 2
